$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 18888.889
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15936
$ws.Range("H23").Value = 18888.889
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15468
$ws.Range("H70").Value = 2816.6
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 2770.75
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 8312.25
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -8852.25
$ws.Range("H73").Value = 2816.6
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 2770.75
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 8312.25
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -10184.25
$ws.Range("H76").Value = 4447592.5
$ws.Range("I76").Value = 4833957
$ws.Range("J76").Value = 4399.5
$ws.Range("K76").Value = 4833957
$ws.Range("L76").Value = 4399.5
$ws.Range("M76").Value = -4833642
$ws.Range("N76").Value = -5029.5
$ws.Range("H79").Value = 4447592.5
$ws.Range("I79").Value = 4833957
$ws.Range("J79").Value = 4399.5
$ws.Range("K79").Value = 4833957
$ws.Range("L79").Value = 4399.5
$ws.Range("M79").Value = -4832865
$ws.Range("N79").Value = -6583.5
$ws.Range("H138").Value = 1848.29
$ws.Range("I138").Value = 621.09375
$ws.Range("J138").Value = 2425.7942
$ws.Range("K138").Value = 1863.28125
$ws.Range("L138").Value = 7277.382599999999
$ws.Range("M138").Value = 3276.71875
$ws.Range("N138").Value = -17557.3826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2000134
$ws.Range("I5").Value = 2000134
$ws.Range("K5").Value = 2000134
$ws.Range("M5").Value = -2000022
$ws.Range("H32").Value = 2586.1973
$ws.Range("I32").Value = 1738.541
$ws.Range("J32").Value = 6033.3335
$ws.Range("K32").Value = 1738.541
$ws.Range("L32").Value = 6033.3335
$ws.Range("M32").Value = -1451.541
$ws.Range("N32").Value = -6607.3335
$ws.Range("H61").Value = 2794.3447
$ws.Range("I61").Value = 1928.5555
$ws.Range("J61").Value = 4211.091
$ws.Range("K61").Value = 1928.5555
$ws.Range("L61").Value = 4211.091
$ws.Range("M61").Value = -1716.5555
$ws.Range("N61").Value = -4635.091
$ws.Range("H110").Value = 780
$ws.Range("I110").Value = 826.6799999999999
$ws.Range("J110").Value = 585.5
$ws.Range("K110").Value = 826.6799999999999
$ws.Range("L110").Value = 585.5
$ws.Range("M110").Value = 1218.32
$ws.Range("N110").Value = -4675.5
$ws.Range("H132").Value = 3787.5293
$ws.Range("I132").Value = 3138.4
$ws.Range("J132").Value = 4714.857
$ws.Range("K132").Value = 9415.200000000001
$ws.Range("L132").Value = 14144.571
$ws.Range("M132").Value = -6885.200000000001
$ws.Range("N132").Value = -19204.571
$ws.Range("H136").Value = 2794.3447
$ws.Range("I136").Value = 1928.5555
$ws.Range("J136").Value = 4211.091
$ws.Range("K136").Value = 5785.666499999999
$ws.Range("L136").Value = 12633.273
$ws.Range("M136").Value = -3235.666499999999
$ws.Range("N136").Value = -17733.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2000134
$ws.Range("I4").Value = 2000134
$ws.Range("K4").Value = 2000134
$ws.Range("M4").Value = -2000019
$ws.Range("H22").Value = 94
$ws.Range("I22").Value = 57.857143
$ws.Range("J22").Value = 178.33333
$ws.Range("K22").Value = 57.857143
$ws.Range("L22").Value = 178.33333
$ws.Range("M22").Value = 115.142857
$ws.Range("N22").Value = -524.3333299999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 6696666.5
$ws.Range("I3").Value = 6696666.5
$ws.Range("K3").Value = 6696666.5
$ws.Range("M3").Value = -6696553.5
$ws.Range("H31").Value = 1234.697
$ws.Range("I31").Value = 911.9796
$ws.Range("J31").Value = 1550.96
$ws.Range("K31").Value = 911.9796
$ws.Range("L31").Value = 1550.96
$ws.Range("M31").Value = -616.9796
$ws.Range("N31").Value = -2140.96
$ws.Range("H34").Value = 1234.697
$ws.Range("I34").Value = 911.9796
$ws.Range("J34").Value = 1550.96
$ws.Range("K34").Value = 911.9796
$ws.Range("L34").Value = 1550.96
$ws.Range("M34").Value = -709.9796
$ws.Range("N34").Value = -1954.96
$ws.Range("H105").Value = 384.10526
$ws.Range("I105").Value = 305.66666
$ws.Range("K105").Value = 305.66666
$ws.Range("M105").Value = 1441.33334
$ws.Range("H132").Value = 2665.25
$ws.Range("I132").Value = 1987.6666
$ws.Range("J132").Value = 3958.818
$ws.Range("K132").Value = 5962.9998
$ws.Range("L132").Value = 11876.454
$ws.Range("M132").Value = -3432.9998
$ws.Range("N132").Value = -16936.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1356.1111
$ws.Range("I5").Value = 593.63635
$ws.Range("J5").Value = 2085.4348
$ws.Range("K5").Value = 1780.90905
$ws.Range("L5").Value = 6256.3044
$ws.Range("M5").Value = -1668.90905
$ws.Range("N5").Value = -6480.3044
$ws.Range("H68").Value = 837.5599999999999
$ws.Range("I68").Value = 723.08826
$ws.Range("J68").Value = 1080.8125
$ws.Range("K68").Value = 2169.26478
$ws.Range("L68").Value = 3242.4375
$ws.Range("M68").Value = -1358.26478
$ws.Range("N68").Value = -4864.4375
$ws.Range("H71").Value = 837.5599999999999
$ws.Range("I71").Value = 723.08826
$ws.Range("J71").Value = 1080.8125
$ws.Range("K71").Value = 6507.79434
$ws.Range("L71").Value = 9727.3125
$ws.Range("M71").Value = -2451.79434
$ws.Range("N71").Value = -17839.3125
$ws.Range("H107").Value = 1012.7544
$ws.Range("I107").Value = 364
$ws.Range("J107").Value = 1244.4524
$ws.Range("K107").Value = 1092
$ws.Range("L107").Value = 3733.357199999999
$ws.Range("M107").Value = 828
$ws.Range("N107").Value = -7573.357199999999
$ws.Range("H116").Value = 1719.8
$ws.Range("J116").Value = 2825
$ws.Range("L116").Value = 8475
$ws.Range("N116").Value = -15359
$ws.Range("H117").Value = 1900
$ws.Range("I117").Value = 400
$ws.Range("J117").Value = 2650
$ws.Range("K117").Value = 1200
$ws.Range("L117").Value = 7950
$ws.Range("M117").Value = 2242
$ws.Range("N117").Value = -14834
$ws.Range("H119").Value = 1800
$ws.Range("I119").Value = 800
$ws.Range("J119").Value = 2800
$ws.Range("K119").Value = 2400
$ws.Range("L119").Value = 8400
$ws.Range("M119").Value = 2438
$ws.Range("N119").Value = -18076
$ws.Range("H120").Value = 35171.668
$ws.Range("I120").Value = 35343.332
$ws.Range("K120").Value = 106029.996
$ws.Range("M120").Value = -101191.996
$ws.Range("H132").Value = 1161.2941
$ws.Range("J132").Value = 1224.2142
$ws.Range("L132").Value = 11017.9278
$ws.Range("N132").Value = -16077.9278
$ws.Range("H135").Value = 1356.1111
$ws.Range("I135").Value = 593.63635
$ws.Range("J135").Value = 2085.4348
$ws.Range("K135").Value = 5342.72715
$ws.Range("L135").Value = 18768.9132
$ws.Range("M135").Value = -2807.72715
$ws.Range("N135").Value = -23838.9132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H70").Value = 6409.85
$ws.Range("I70").Value = 7314.143
$ws.Range("J70").Value = 4299.8335
$ws.Range("K70").Value = 7314.143
$ws.Range("L70").Value = 4299.8335
$ws.Range("M70").Value = -7044.143
$ws.Range("N70").Value = -4839.8335
$ws.Range("H73").Value = 6409.85
$ws.Range("I73").Value = 7314.143
$ws.Range("J73").Value = 4299.8335
$ws.Range("K73").Value = 7314.143
$ws.Range("L73").Value = 4299.8335
$ws.Range("M73").Value = -6378.143
$ws.Range("N73").Value = -6171.8335
$ws.Range("H122").Value = 11111111
$ws.Range("I122").Value = 11111111
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 33333333
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -33330883
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 919.2222
$ws.Range("I22").Value = 981.8333
$ws.Range("J22").Value = 794
$ws.Range("K22").Value = 981.8333
$ws.Range("L22").Value = 794
$ws.Range("M22").Value = -686.8333
$ws.Range("N22").Value = -1384
$ws.Range("H27").Value = 919.2222
$ws.Range("I27").Value = 981.8333
$ws.Range("J27").Value = 794
$ws.Range("K27").Value = 981.8333
$ws.Range("L27").Value = 794
$ws.Range("M27").Value = -874.8333
$ws.Range("N27").Value = -1008
$ws.Range("H46").Value = 497.6154
$ws.Range("I46").Value = 514
$ws.Range("J46").Value = 487.375
$ws.Range("K46").Value = 514
$ws.Range("L46").Value = 487.375
$ws.Range("M46").Value = -326
$ws.Range("N46").Value = -863.375
$ws.Range("H55").Value = 441.1111
$ws.Range("I55").Value = 467
$ws.Range("K55").Value = 467
$ws.Range("M55").Value = -294
$ws.Range("H122").Value = 3884.6155
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3884.6155
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11653.8465
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -16553.8465
$ws.Range("H132").Value = 3381.5806
$ws.Range("I132").Value = 2395.6956
$ws.Range("J132").Value = 6216
$ws.Range("K132").Value = 7187.0868
$ws.Range("L132").Value = 18648
$ws.Range("M132").Value = -4657.0868
$ws.Range("N132").Value = -23708

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 29996.334
$ws.Range("J20").Value = 29996.334
$ws.Range("L20").Value = 29996.334
$ws.Range("N20").Value = -30476.334
$ws.Range("H22").Value = 30000
$ws.Range("J22").Value = 30000
$ws.Range("L22").Value = 30000
$ws.Range("N22").Value = -30586
$ws.Range("H61").Value = 5878.7144
$ws.Range("I61").Value = 3850.3333
$ws.Range("K61").Value = 3850.3333
$ws.Range("M61").Value = -3558.3333
$ws.Range("H132").Value = 15627010
$ws.Range("I132").Value = 18519982
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 55559946
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -55557416
$ws.Range("N132").Value = -19964
$ws.Range("H136").Value = 15922108
$ws.Range("I136").Value = 19667342
$ws.Range("J136").Value = 4863.25
$ws.Range("K136").Value = 59002026
$ws.Range("L136").Value = 14589.75
$ws.Range("M136").Value = -58999476
$ws.Range("N136").Value = -19689.75
